$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.084.36"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "1.825.12"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").Value = "'312.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "

$ws.Range("D6").Value = "'1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").Value = "'0.4683"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").Value = "'0.3651"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").Value = "'0.07383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("D10").Value = "'0.8786"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").Value = "'20.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "1.871.29"
$ws.Range("E12").Value = "  +2.28%  "

$ws.Range("D13").Value = "'0.07437"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.07%  "

$ws.Range("D14").Value = "'5.368"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.91%  "

$ws.Range("D15").Value = "'92.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").Value = "'6.521"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "'0.000008711"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").Value = "27.566.02"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("D21").Value = "'14.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("D22").Value = "'5.234"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").Value = "'10.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "2.080.84"
$ws.Range("E24").Value = "  +1.38%  "

$ws.Range("D25").Value = "'1.882"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("D26").Value = "'151.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").Value = "'18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").Value = "'2.136"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").Value = "'5.154"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.61%  "

$ws.Range("D30").Value = "'116.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.59%  "

$ws.Range("D31").Value = "'0.08887"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

$ws.Range("D32").Value = "'0.7433"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").Value = "'1.162"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").Value = "'4.504"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("D35").Value = "'2.940"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").Value = "'1.007"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").Value = "'2.530"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.92%  "

$ws.Range("D38").Value = "'1.089"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.57%  "

$ws.Range("D39").Value = "'0.05288"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.36%  "

$ws.Range("D40").Value = "'0.01931"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "'7.329"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").Value = "'2.931"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("D43").Value = "'0.5249"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "'0.1638"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.76%  "

$ws.Range("D45").Value = "'8.360"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").Value = "'0.4893"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "'10.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").Value = "'1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "'104.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("D50").Value = "'1.649"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").Value = "'0.06267"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "

